# 19/12/2025: Update the list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-30 (A:E) representing the refreshed pipeline list.
$data = @(
    @(196, "Blockaid", "Enterprise Account Executive (Fintech)", "Jacob Neuhard", "CV Sent"),
    @(416, "Prime Security", "Account Executive", "David Chase", "CV Sent"),
    @(484, "Cognition AI", "Enterprise Account Executive (US)", "Rich Kroll", "CV Sent"),
    @(541, "Blockaid", "CSM - Singapore", "Adeline  Wong", "CV Sent"),
    @(572, "Chainguard", "Strategic Account Manager (Tola)", "Rich Kroll", "CV Sent"),
    @(572, "Chainguard", "Strategic Account Manager (Tola)", "Todd Ferreira", "CV Sent"),
    @(577, "Oasis Security", "TAM (East + Central)", "Aaron Lake", "1st Interview"),
    @(617, "Dash0", "Commercial AE NYC", "Daniel Mehta", "CV Sent"),
    @(755, "Vega Security", "Enterprise AE (Mid-ATL)", "Rich Kroll", "CV Sent"),
    @(761, "Mate Security", "VP Sales - US", "Amanda Parks", "1st Interview"),
    @(761, "Mate Security", "VP Sales - US", "Andrew Kashin", "CV Sent"),
    @(761, "Mate Security", "VP Sales - US", "Brad Rhodes", "CV Sent"),
    @(761, "Mate Security", "VP Sales - US", "Jeff Ross", "1st Interview"),
    @(769, "TraceBit", "VP Sales NYC", "Jourdan Loffredo", "1st Interview"),
    @(769, "TraceBit", "VP Sales NYC", "Eric Amstutz", "1st Interview"),
    @(770, "TraceBit", "Account Executive NYC", "Jacob Neuhard", "CV Sent"),
    @(770, "TraceBit", "Account Executive NYC", "Matthew Bohannon", "1st Interview"),
    @(770, "TraceBit", "Account Executive NYC", "Ryan Drillock", "3rd Interview"),
    @(791, "Adaptive6", "Head of Sales (US)", "Andrew Kashin", "CV Sent"),
    @(809, "Blink Ops", "Enterprise Account Executive (Mid-ATL)", "Rich Kroll", "CV Sent"),
    @(834, "Blockaid", "Regional Director US", "Chris Dyer", "CV Sent"),
    @(834, "Blockaid", "Regional Director US", "Jeff White", "1st Interview"),
    @(841, "Maze", "Technical Success Manager", "Aaron Lake", "2nd Interview"),
    @(841, "Maze", "Technical Success Manager", "Christian Gerling", "CV Sent"),
    @(841, "Maze", "Technical Success Manager", "James Sandwick", "1st Interview"),
    @(841, "Maze", "Technical Success Manager", "Vishal Soni", "1st Interview"),
    @(843, "Red Access", "Enterprise Account Executive UK", "Michele Maioli", "1st Interview"),
    @(847, "Simile.ai", "CS1 Simile.ai - Enterprise AE x2", "Colby Perry", "CV Sent"),
    @(870, "Native Security", "Enterprise Account Executive (North-East)", "Rich Kroll", "CV Sent")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}

# The table shrank from 32 to 30 rows (A1:E32 -> A1:E30); clear the two
# trailing rows that no longer hold data.
$ws.Range("A31:E32").Clear()
